# "incluindo novo arquivo de validacao"
#
# The workbook currently has 3 tabs, in this order:
#   1. DadosCarga   (Nome,CPF,Data_Nascimento,Sexo,Matricula)                 order: Mario, Anthony, Yasmin, Marina
#   2. DadosSistema (same 5 cols + chaveCPF&Nome)                             order: Anthony, Mario, Yasmin, Marina
#   3. Analise      (same 5 cols + ComparaCPF, chaveCPF&Nome, ComparaCPF_Nome) order: Mario, Anthony, Yasmin, Marina
#
# After the edit the tabs become (same physical order, first tab renamed to
# front): Analise, DadosCarga, DadosSistema -- each sheet's underlying data
# also changes:
#   1. Analise      -> the old DadosCarga data, PLUS 11 new comparison/key
#                      columns (F..P), each row "Localizado".
#   2. DadosCarga   -> the old DadosSistema's 5 base columns but with rows
#                      re-ordered to Mario, Anthony, Yasmin, Marina and the
#                      old chaveCPF&Nome column (F) removed.
#   3. DadosSistema -> the old Analise's 5 base columns, re-ordered back to
#                      Anthony, Mario, Yasmin, Marina, with 5 new "chave*"
#                      columns (F..J) replacing the old Compara*/chave cols.

$wb = $excel.ActiveWorkbook

$wsCarga    = $wb.Worksheets.Item("DadosCarga")
$wsSistema  = $wb.Worksheets.Item("DadosSistema")
$wsAnalise  = $wb.Worksheets.Item("Analise")

# --- rotate the three sheet names (use scratch names so we never collide) ---
$wsCarga.Name   = "__rot_carga__"
$wsSistema.Name = "__rot_sistema__"
$wsAnalise.Name = "__rot_analise__"

$wsCarga.Name   = "Analise"       # was DadosCarga
$wsSistema.Name = "DadosCarga"    # was DadosSistema
$wsAnalise.Name = "DadosSistema"  # was Analise

$newAnalise  = $wb.Worksheets.Item("Analise")
$newCarga    = $wb.Worksheets.Item("DadosCarga")
$newSistema  = $wb.Worksheets.Item("DadosSistema")

# The four people, keyed by name, with every derived value needed below.
$mario = @{
    Nome = "Mário Manoel Calebe Moura"; CPF = "039.173.221-88"; Dt = 19993
    DtStr = "1954-09-26"; Sexo = "M"; Mat = 70297
}
$anthony = @{
    Nome = "Anthony Henrique Costa"; CPF = "876.995.634-09"; Dt = 33689
    DtStr = "1992-03-26"; Sexo = "M"; Mat = 69314
}
$yasmin = @{
    Nome = "Yasmin Eliane Agatha Gonçalves"; CPF = "193.703.911-00"; Dt = 16635
    DtStr = "1945-07-17"; Sexo = "F"; Mat = 61902
}
$marina = @{
    Nome = "Marina Sophie Marlene da Luz"; CPF = "369.960.476-41"; Dt = 16107
    DtStr = "1944-02-05"; Sexo = "F"; Mat = 38045
}

function Set-BaseRow($ws, $row, $p) {
    $ws.Cells.Item($row, 1).Value = $p.Nome
    $ws.Cells.Item($row, 2).Value = $p.CPF
    $ws.Cells.Item($row, 3).Value = $p.Dt
    $ws.Cells.Item($row, 4).Value = $p.Sexo
    $ws.Cells.Item($row, 5).Value = $p.Mat
}

# =====================================================================
# 1) "Analise" (new) -- old DadosCarga rows (Mario,Anthony,Yasmin,Marina)
#    already sit in A:E; add the 11 comparison/key columns F..P.
# =====================================================================
$order1 = @($mario, $anthony, $yasmin, $marina)

$newAnalise.Cells.Item(1, 6).Value  = "ComparaCPF"
$newAnalise.Cells.Item(1, 7).Value  = "chaveCPF&Nome"
$newAnalise.Cells.Item(1, 8).Value  = "ComparaCPF_Nome"
$newAnalise.Cells.Item(1, 9).Value  = "chaveCPF&DtNasc"
$newAnalise.Cells.Item(1, 10).Value = "ComparaCPF_DtNasc"
$newAnalise.Cells.Item(1, 11).Value = "chaveCPF&Sexo"
$newAnalise.Cells.Item(1, 12).Value = "ComparaCPF_Sexo"
$newAnalise.Cells.Item(1, 13).Value = "chaveCPF&Matricula"
$newAnalise.Cells.Item(1, 14).Value = "ComparaCPF_Matricula"
$newAnalise.Cells.Item(1, 15).Value = "chave"
$newAnalise.Cells.Item(1, 16).Value = "Resultado"
$newAnalise.Cells.Item(1, 5).Copy()
$newAnalise.Range("F1:P1").PasteSpecial(-4122)

for ($i = 0; $i -lt $order1.Length; $i++) {
    $p = $order1[$i]
    $row = $i + 2
    Set-BaseRow $newAnalise $row $p

    $chaveNome = $p.CPF + $p.Nome
    $chaveDt   = $p.CPF + $p.DtStr
    $chaveSexo = $p.CPF + $p.Sexo
    $chaveMat  = $p.CPF + [string]$p.Mat
    $chave     = $p.CPF + $p.Nome + $p.DtStr + $p.Sexo + [string]$p.Mat

    $newAnalise.Cells.Item($row, 6).Value  = "Localizado"
    $newAnalise.Cells.Item($row, 7).Value  = $chaveNome
    $newAnalise.Cells.Item($row, 8).Value  = "Localizado"
    $newAnalise.Cells.Item($row, 9).Value  = $chaveDt
    $newAnalise.Cells.Item($row, 10).Value = "Localizado"
    $newAnalise.Cells.Item($row, 11).Value = $chaveSexo
    $newAnalise.Cells.Item($row, 12).Value = "Localizado"
    $newAnalise.Cells.Item($row, 13).Value = $chaveMat
    $newAnalise.Cells.Item($row, 14).Value = "Localizado"
    $newAnalise.Cells.Item($row, 15).Value = $chave
    $newAnalise.Cells.Item($row, 16).Value = "Localizado"
}

# =====================================================================
# 2) "DadosCarga" (new) -- old DadosSistema data: re-order rows back to
#    Mario,Anthony,Yasmin,Marina and drop the old chaveCPF&Nome column.
# =====================================================================
$order2 = @($mario, $anthony, $yasmin, $marina)
for ($i = 0; $i -lt $order2.Length; $i++) {
    $row = $i + 2
    Set-BaseRow $newCarga $row $order2[$i]
}
$newCarga.Columns.Item(6).Delete()

# =====================================================================
# 3) "DadosSistema" (new) -- old Analise data: re-order rows to
#    Anthony,Mario,Yasmin,Marina and replace F:H with 5 new chave* columns.
# =====================================================================
$newSistema.Columns.Item(8).Delete()
$newSistema.Columns.Item(7).Delete()
$newSistema.Columns.Item(6).Delete()

$order3 = @($anthony, $mario, $yasmin, $marina)

$newSistema.Cells.Item(1, 6).Value  = "chaveCPF&Nome"
$newSistema.Cells.Item(1, 7).Value  = "chaveCPF&DtNasc"
$newSistema.Cells.Item(1, 8).Value  = "chaveCPF&Sexo"
$newSistema.Cells.Item(1, 9).Value  = "chaveCPF&Matricula"
$newSistema.Cells.Item(1, 10).Value = "chave"
$newSistema.Cells.Item(1, 5).Copy()
$newSistema.Range("F1:J1").PasteSpecial(-4122)

for ($i = 0; $i -lt $order3.Length; $i++) {
    $p = $order3[$i]
    $row = $i + 2
    Set-BaseRow $newSistema $row $p

    $chaveNome = $p.CPF + $p.Nome
    $chaveDt   = $p.CPF + $p.DtStr
    $chaveSexo = $p.CPF + $p.Sexo
    $chaveMat  = $p.CPF + [string]$p.Mat
    $chave     = $p.CPF + $p.Nome + $p.DtStr + $p.Sexo + [string]$p.Mat

    $newSistema.Cells.Item($row, 6).Value  = $chaveNome
    $newSistema.Cells.Item($row, 7).Value  = $chaveDt
    $newSistema.Cells.Item($row, 8).Value  = $chaveSexo
    $newSistema.Cells.Item($row, 9).Value  = $chaveMat
    $newSistema.Cells.Item($row, 10).Value = $chave
}
